# The pre-bound $wb / $ws handles are not reliably populated in this
# runtime, so re-derive them from the Application object.
$wb = $excel.ActiveWorkbook

# Locate the data sheet by its current name (second tab, after the blank "Sheet").
$ws = $wb.Worksheets.Item("scenario_3_RL")

# Rename the sheet to match the new scenario label.
$ws.Name = "scenario_1_RL"

# Helper: write a numeric-looking value while keeping the cell a TEXT cell
# (matches the source file's inline-string / text-typed cells), and strip
# the Text number-format we had to apply so no stray style sticks around.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("B1") "62.25"
Set-TextValue $ws.Range("B2") "1.19"
Set-TextValue $ws.Range("B3") "3.74"
Set-TextValue $ws.Range("B4") "4"
